$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D
$ws.Range("D1").Value = "ITI"

# New ConditionType (column C) values for rows 2-17
$conditionType = @(1,4,4,3,3,1,2,1,4,2,2,4,1,3,2,3)
# New ITI (column D) values for rows 2-17
$iti = @(7,7,6,9,6,6,9,6,7,7,7,7,8,6,6,8)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $conditionType[$i]
    $ws.Cells.Item($row, 4).Value = $iti[$i]
}

# Remove rows 18, 19, 20 (trailing rows no longer present)
$ws.Range("A18:D20").Clear()

$ws.Range("D18").Select()
